# Auto-generated Excel COM-interop script to apply scheduled market-data
# refresh updates to the per-sheet "Leve profit" tables (columns H-N).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 15626817   # H70: 15626892 -> 15626817
$ws.Cells.Item(70, 9).Value = 2246.8333   # I70: 2538 -> 2246.8333
$ws.Cells.Item(70, 10).Value = 25001558   # J70: 22728870 -> 25001558
$ws.Cells.Item(70, 11).Value = 6740.499899999999   # K70: 7614 -> 6740.499899999999
$ws.Cells.Item(70, 12).Value = 75004674   # L70: 68186610 -> 75004674
$ws.Cells.Item(70, 13).Value = -6470.499899999999   # M70: -7344 -> -6470.499899999999
$ws.Cells.Item(70, 14).Value = -75005214   # N70: -68187150 -> -75005214
$ws.Cells.Item(73, 8).Value = 15626817   # H73: 15626892 -> 15626817
$ws.Cells.Item(73, 9).Value = 2246.8333   # I73: 2538 -> 2246.8333
$ws.Cells.Item(73, 10).Value = 25001558   # J73: 22728870 -> 25001558
$ws.Cells.Item(73, 11).Value = 6740.499899999999   # K73: 7614 -> 6740.499899999999
$ws.Cells.Item(73, 12).Value = 75004674   # L73: 68186610 -> 75004674
$ws.Cells.Item(73, 13).Value = -5804.499899999999   # M73: -6678 -> -5804.499899999999
$ws.Cells.Item(73, 14).Value = -75006546   # N73: -68188482 -> -75006546
$ws.Cells.Item(98, 8).Value = 1222.7241   # H98: 1262.4642 -> 1222.7241
$ws.Cells.Item(98, 9).Value = 1120.7037   # I98: 1159.5769 -> 1120.7037
$ws.Cells.Item(98, 11).Value = 1120.7037   # K98: 1159.5769 -> 1120.7037
$ws.Cells.Item(98, 13).Value = 377.2963   # M98: 338.4231 -> 377.2963
$ws.Cells.Item(111, 8).Value = 1996.3334   # H111: 2150.3333 -> 1996.3334
$ws.Cells.Item(111, 9).Value = 1996.3334   # I111: 2262 -> 1996.3334
$ws.Cells.Item(111, 10).Value = 0   # J111: 1927 -> 0
$ws.Cells.Item(111, 11).Value = 5989.0002   # K111: 6786 -> 5989.0002
$ws.Cells.Item(111, 12).Value = 0   # L111: 5781 -> 0
$ws.Cells.Item(111, 13).Value = -2922.0002   # M111: -3719 -> -2922.0002
$ws.Cells.Item(111, 14).ClearContents()   # N111: delete (was -11915)
$ws.Cells.Item(113, 8).Value = 2589.7273   # H113: 2399 -> 2589.7273
$ws.Cells.Item(113, 9).Value = 2272.125   # I113: 2159.5386 -> 2272.125
$ws.Cells.Item(113, 11).Value = 2272.125   # K113: 2159.5386 -> 2272.125
$ws.Cells.Item(113, 13).Value = 981.875   # M113: 1094.4614 -> 981.875
$ws.Cells.Item(122, 8).Value = 1222.7241   # H122: 1262.4642 -> 1222.7241
$ws.Cells.Item(122, 9).Value = 1120.7037   # I122: 1159.5769 -> 1120.7037
$ws.Cells.Item(122, 11).Value = 3362.1111   # K122: 3478.7307 -> 3362.1111
$ws.Cells.Item(122, 13).Value = -912.1111000000001   # M122: -1028.7307 -> -912.1111000000001
$ws.Cells.Item(138, 8).Value = 1826.1   # H138: 1906.99 -> 1826.1
$ws.Cells.Item(138, 9).Value = 1060.2667   # I138: 1100.0238 -> 1060.2667
$ws.Cells.Item(138, 10).Value = 2452.691   # J138: 2491.3447 -> 2452.691
$ws.Cells.Item(138, 11).Value = 3180.800099999999   # K138: 3300.0714 -> 3180.800099999999
$ws.Cells.Item(138, 12).Value = 7358.072999999999   # L138: 7474.034100000001 -> 7358.072999999999
$ws.Cells.Item(138, 13).Value = 1959.199900000001   # M138: 1839.9286 -> 1959.199900000001
$ws.Cells.Item(138, 14).Value = -17638.073   # N138: -17754.0341 -> -17638.073

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1732.7916   # H2: 1742.7916 -> 1732.7916
$ws.Cells.Item(2, 9).Value = 1146.6842   # I2: 1159.3158 -> 1146.6842
$ws.Cells.Item(2, 11).Value = 1146.6842   # K2: 1159.3158 -> 1146.6842
$ws.Cells.Item(2, 13).Value = -1033.6842   # M2: -1046.3158 -> -1033.6842
$ws.Cells.Item(32, 8).Value = 18660252   # H32: 18660378 -> 18660252
$ws.Cells.Item(32, 9).Value = 18384086   # I32: 19064966 -> 18384086
$ws.Cells.Item(32, 10).Value = 23815356   # J32: 14290813 -> 23815356
$ws.Cells.Item(32, 11).Value = 18384086   # K32: 19064966 -> 18384086
$ws.Cells.Item(32, 12).Value = 23815356   # L32: 14290813 -> 23815356
$ws.Cells.Item(32, 13).Value = -18383799   # M32: -19064679 -> -18383799
$ws.Cells.Item(32, 14).Value = -23815930   # N32: -14291387 -> -23815930
$ws.Cells.Item(74, 8).Value = 2290.9756   # H74: 2330.725 -> 2290.9756
$ws.Cells.Item(74, 9).Value = 2364.353   # I74: 2416.2727 -> 2364.353
$ws.Cells.Item(74, 10).Value = 1934.5714   # J74: 1927.4286 -> 1934.5714
$ws.Cells.Item(74, 11).Value = 2364.353   # K74: 2416.2727 -> 2364.353
$ws.Cells.Item(74, 12).Value = 1934.5714   # L74: 1927.4286 -> 1934.5714
$ws.Cells.Item(74, 13).Value = -1490.353   # M74: -1542.2727 -> -1490.353
$ws.Cells.Item(74, 14).Value = -3682.5714   # N74: -3675.4286 -> -3682.5714
$ws.Cells.Item(77, 8).Value = 2290.9756   # H77: 2330.725 -> 2290.9756
$ws.Cells.Item(77, 9).Value = 2364.353   # I77: 2416.2727 -> 2364.353
$ws.Cells.Item(77, 10).Value = 1934.5714   # J77: 1927.4286 -> 1934.5714
$ws.Cells.Item(77, 11).Value = 11821.765   # K77: 12081.3635 -> 11821.765
$ws.Cells.Item(77, 12).Value = 9672.857   # L77: 9637.143 -> 9672.857
$ws.Cells.Item(77, 13).Value = -7453.764999999999   # M77: -7713.363499999999 -> -7453.764999999999
$ws.Cells.Item(77, 14).Value = -18408.857   # N77: -18373.143 -> -18408.857
$ws.Cells.Item(102, 8).Value = 1601   # H102: 1695.95 -> 1601
$ws.Cells.Item(116, 8).Value = 1732.7916   # H116: 1742.7916 -> 1732.7916
$ws.Cells.Item(116, 9).Value = 1146.6842   # I116: 1159.3158 -> 1146.6842
$ws.Cells.Item(116, 11).Value = 1146.6842   # K116: 1159.3158 -> 1146.6842
$ws.Cells.Item(116, 13).Value = 1147.3158   # M116: 1134.6842 -> 1147.3158
$ws.Cells.Item(122, 8).Value = 3871.4055   # H122: 4013.6943 -> 3871.4055
$ws.Cells.Item(122, 9).Value = 2181.1904   # I122: 2252.75 -> 2181.1904
$ws.Cells.Item(122, 10).Value = 6089.8125   # J122: 6214.875 -> 6089.8125
$ws.Cells.Item(122, 11).Value = 6543.5712   # K122: 6758.25 -> 6543.5712
$ws.Cells.Item(122, 12).Value = 18269.4375   # L122: 18644.625 -> 18269.4375
$ws.Cells.Item(122, 13).Value = -4093.5712   # M122: -4308.25 -> -4093.5712
$ws.Cells.Item(122, 14).Value = -23169.4375   # N122: -23544.625 -> -23169.4375
$ws.Cells.Item(130, 8).Value = 39295.09   # H130: 36359.125 -> 39295.09
$ws.Cells.Item(130, 10).Value = 39295.09   # J130: 36359.125 -> 39295.09
$ws.Cells.Item(130, 12).Value = 39295.09   # L130: 36359.125 -> 39295.09
$ws.Cells.Item(130, 14).Value = -49335.09   # N130: -46399.125 -> -49335.09
$ws.Cells.Item(135, 8).Value = 0   # H135: 99900 -> 0
$ws.Cells.Item(135, 10).Value = 0   # J135: 99900 -> 0
$ws.Cells.Item(135, 12).Value = 0   # L135: 99900 -> 0
$ws.Cells.Item(135, 14).ClearContents()   # N135: delete (was -110040)

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1732.7916   # H3: 1742.7916 -> 1732.7916
$ws.Cells.Item(3, 9).Value = 1146.6842   # I3: 1159.3158 -> 1146.6842
$ws.Cells.Item(3, 11).Value = 1146.6842   # K3: 1159.3158 -> 1146.6842
$ws.Cells.Item(3, 13).Value = -1032.6842   # M3: -1045.3158 -> -1032.6842

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1046.5   # H16: 1089.9412 -> 1046.5
$ws.Cells.Item(16, 9).Value = 852.26666   # I16: 891.1429000000001 -> 852.26666
$ws.Cells.Item(16, 11).Value = 852.26666   # K16: 891.1429000000001 -> 852.26666
$ws.Cells.Item(16, 13).Value = -565.26666   # M16: -604.1429000000001 -> -565.26666
$ws.Cells.Item(58, 8).Value = 2739.3438   # H58: 2838.6333 -> 2739.3438
$ws.Cells.Item(58, 9).Value = 2691.074   # I58: 2864.9583 -> 2691.074
$ws.Cells.Item(58, 10).Value = 3000   # J58: 2733.3333 -> 3000
$ws.Cells.Item(58, 11).Value = 2691.074   # K58: 2864.9583 -> 2691.074
$ws.Cells.Item(58, 12).Value = 3000   # L58: 2733.3333 -> 3000
$ws.Cells.Item(58, 13).Value = -2488.074   # M58: -2661.9583 -> -2488.074
$ws.Cells.Item(58, 14).Value = -3406   # N58: -3139.3333 -> -3406
$ws.Cells.Item(105, 8).Value = 1696.1111   # H105: 1749.5294 -> 1696.1111
$ws.Cells.Item(105, 9).Value = 1232.3077   # I105: 1269.3334 -> 1232.3077
$ws.Cells.Item(105, 11).Value = 1232.3077   # K105: 1269.3334 -> 1232.3077
$ws.Cells.Item(105, 13).Value = 514.6922999999999   # M105: 477.6666 -> 514.6922999999999
$ws.Cells.Item(113, 8).Value = 1046.5   # H113: 1089.9412 -> 1046.5
$ws.Cells.Item(113, 9).Value = 852.26666   # I113: 891.1429000000001 -> 852.26666
$ws.Cells.Item(113, 11).Value = 852.26666   # K113: 891.1429000000001 -> 852.26666
$ws.Cells.Item(113, 13).Value = 1317.73334   # M113: 1278.8571 -> 1317.73334
$ws.Cells.Item(122, 8).Value = 3849211   # H122: 3574286.8 -> 3849211
$ws.Cells.Item(122, 9).Value = 4765121   # I122: 4169545 -> 4765121
$ws.Cells.Item(122, 10).Value = 2389.4   # J122: 2737.25 -> 2389.4
$ws.Cells.Item(122, 11).Value = 14295363   # K122: 12508635 -> 14295363
$ws.Cells.Item(122, 12).Value = 7168.200000000001   # L122: 8211.75 -> 7168.200000000001
$ws.Cells.Item(122, 13).Value = -14292913   # M122: -12506185 -> -14292913
$ws.Cells.Item(122, 14).Value = -12068.2   # N122: -13111.75 -> -12068.2
$ws.Cells.Item(132, 8).Value = 2694.7673   # H132: 2745.3572 -> 2694.7673
$ws.Cells.Item(132, 9).Value = 2624.103   # I132: 2686.3484 -> 2624.103
$ws.Cells.Item(132, 11).Value = 7872.309   # K132: 8059.0452 -> 7872.309
$ws.Cells.Item(132, 13).Value = -5342.309   # M132: -5529.0452 -> -5342.309
$ws.Cells.Item(134, 8).Value = 2157.0264   # H134: 2158.0789 -> 2157.0264
$ws.Cells.Item(134, 9).Value = 1893.3939   # I134: 1894.6061 -> 1893.3939
$ws.Cells.Item(134, 11).Value = 5680.1817   # K134: 5683.8183 -> 5680.1817
$ws.Cells.Item(134, 13).Value = -3145.1817   # M134: -3148.8183 -> -3145.1817
$ws.Cells.Item(136, 8).Value = 2739.3438   # H136: 2838.6333 -> 2739.3438
$ws.Cells.Item(136, 9).Value = 2691.074   # I136: 2864.9583 -> 2691.074
$ws.Cells.Item(136, 10).Value = 3000   # J136: 2733.3333 -> 3000
$ws.Cells.Item(136, 11).Value = 8073.222   # K136: 8594.874899999999 -> 8073.222
$ws.Cells.Item(136, 12).Value = 9000   # L136: 8199.999899999999 -> 9000
$ws.Cells.Item(136, 13).Value = -5523.222   # M136: -6044.874899999999 -> -5523.222
$ws.Cells.Item(136, 14).Value = -14100   # N136: -13299.9999 -> -14100
$ws.Cells.Item(141, 8).Value = 535711.9   # H141: 344361.6 -> 535711.9
$ws.Cells.Item(141, 10).Value = 570405.5600000001   # J141: 357159.9 -> 570405.5600000001
$ws.Cells.Item(141, 12).Value = 570405.5600000001   # L141: 357159.9 -> 570405.5600000001
$ws.Cells.Item(141, 14).Value = -580765.5600000001   # N141: -367519.9 -> -580765.5600000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 336.23077   # H23: 362.16666 -> 336.23077
$ws.Cells.Item(23, 10).Value = 344.2857   # J23: 397.5 -> 344.2857
$ws.Cells.Item(23, 12).Value = 1032.8571   # L23: 1192.5 -> 1032.8571
$ws.Cells.Item(23, 14).Value = -1502.8571   # N23: -1662.5 -> -1502.8571
$ws.Cells.Item(34, 8).Value = 929.7143   # H34: 1056.7778 -> 929.7143
$ws.Cells.Item(34, 10).Value = 2001   # J34: 1801.2 -> 2001
$ws.Cells.Item(34, 12).Value = 6003   # L34: 5403.6 -> 6003
$ws.Cells.Item(34, 14).Value = -6171   # N34: -5571.6 -> -6171
$ws.Cells.Item(39, 8).Value = 4345.467   # H39: 4421.0625 -> 4345.467
$ws.Cells.Item(39, 10).Value = 4345.467   # J39: 4421.0625 -> 4345.467
$ws.Cells.Item(39, 12).Value = 13036.401   # L39: 13263.1875 -> 13036.401
$ws.Cells.Item(39, 14).Value = -13624.401   # N39: -13851.1875 -> -13624.401
$ws.Cells.Item(55, 8).Value = 2066.3845   # H55: 2025.9286 -> 2066.3845
$ws.Cells.Item(55, 10).Value = 3685.6667   # J55: 3373.4285 -> 3685.6667
$ws.Cells.Item(55, 12).Value = 11057.0001   # L55: 10120.2855 -> 11057.0001
$ws.Cells.Item(55, 14).Value = -11411.0001   # N55: -10474.2855 -> -11411.0001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1110.6842   # H102: 1199.2941 -> 1110.6842
$ws.Cells.Item(102, 9).Value = 1058.1765   # I102: 1151.6 -> 1058.1765
$ws.Cells.Item(102, 11).Value = 1058.1765   # K102: 1151.6 -> 1058.1765
$ws.Cells.Item(102, 13).Value = 563.8235   # M102: 470.4000000000001 -> 563.8235

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 13896382   # H40: 13896449 -> 13896382
$ws.Cells.Item(40, 9).Value = 16673459   # I40: 16673539 -> 16673459
$ws.Cells.Item(40, 11).Value = 16673459   # K40: 16673539 -> 16673459
$ws.Cells.Item(40, 13).Value = -16673323   # M40: -16673403 -> -16673323
$ws.Cells.Item(100, 8).Value = 2602.2   # H100: 2201.8333 -> 2602.2
$ws.Cells.Item(100, 9).Value = 1001   # I100: 800.75 -> 1001
$ws.Cells.Item(100, 11).Value = 1001   # K100: 800.75 -> 1001
$ws.Cells.Item(100, 13).Value = -460   # M100: -259.75 -> -460
$ws.Cells.Item(123, 8).Value = 29888   # H123: 0 -> 29888
$ws.Cells.Item(123, 10).Value = 29888   # J123: 0 -> 29888
$ws.Cells.Item(123, 12).Value = 29888   # L123: 0 -> 29888
$ws.Cells.Item(123, 14).Value = -39688   # N123: new cell (was empty)
$ws.Cells.Item(128, 8).Value = 108969.5   # H128: 109316.336 -> 108969.5
$ws.Cells.Item(128, 10).Value = 108969.5   # J128: 109316.336 -> 108969.5
$ws.Cells.Item(128, 12).Value = 108969.5   # L128: 109316.336 -> 108969.5
$ws.Cells.Item(128, 14).Value = -118929.5   # N128: -119276.336 -> -118929.5
$ws.Cells.Item(132, 8).Value = 104648.45   # H132: 115753.5 -> 104648.45
$ws.Cells.Item(132, 9).Value = 129081.5   # I132: 129222.75 -> 129081.5
$ws.Cells.Item(132, 10).Value = 6916.25   # J132: 7999.5 -> 6916.25
$ws.Cells.Item(132, 11).Value = 387244.5   # K132: 387668.25 -> 387244.5
$ws.Cells.Item(132, 12).Value = 20748.75   # L132: 23998.5 -> 20748.75
$ws.Cells.Item(132, 13).Value = -384714.5   # M132: -385138.25 -> -384714.5
$ws.Cells.Item(132, 14).Value = -25808.75   # N132: -29058.5 -> -25808.75
$ws.Cells.Item(136, 8).Value = 2541.44   # H136: 2730 -> 2541.44
$ws.Cells.Item(136, 9).Value = 2239.4   # I136: 2360.647 -> 2239.4
$ws.Cells.Item(136, 10).Value = 3749.6   # J136: 4299.75 -> 3749.6
$ws.Cells.Item(136, 11).Value = 6718.200000000001   # K136: 7081.941 -> 6718.200000000001
$ws.Cells.Item(136, 12).Value = 11248.8   # L136: 12899.25 -> 11248.8
$ws.Cells.Item(136, 13).Value = -4168.200000000001   # M136: -4531.941 -> -4168.200000000001
$ws.Cells.Item(136, 14).Value = -16348.8   # N136: -17999.25 -> -16348.8
$ws.Cells.Item(140, 8).Value = 334998.5   # H140: 419998 -> 334998.5
$ws.Cells.Item(140, 10).Value = 334998.5   # J140: 419998 -> 334998.5
$ws.Cells.Item(140, 12).Value = 334998.5   # L140: 419998 -> 334998.5
$ws.Cells.Item(140, 14).Value = -345358.5   # N140: -430358 -> -345358.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3098.4375   # H132: 3198.6 -> 3098.4375
$ws.Cells.Item(132, 9).Value = 2826.8215   # I132: 2921.5 -> 2826.8215
$ws.Cells.Item(132, 11).Value = 8480.4645   # K132: 8764.5 -> 8480.4645
$ws.Cells.Item(132, 13).Value = -5950.4645   # M132: -6234.5 -> -5950.4645
